$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NACE 64 labels: prefix each revenue-code row with its NACE section letter
$ws.Range("A2").Value = "A01"
$ws.Range("A3").Value = "A02"
$ws.Range("A4").Value = "A03"
$ws.Range("A5").Value = "B05-09"
$ws.Range("A6").Value = "C10-12"
$ws.Range("A7").Value = "C13-15"
$ws.Range("A8").Value = "C16"
$ws.Range("A9").Value = "C17"
$ws.Range("A10").Value = "C18"
$ws.Range("A11").Value = "C19"
$ws.Range("A12").Value = "C20"
$ws.Range("A13").Value = "C21"
$ws.Range("A14").Value = "C22"
$ws.Range("A15").Value = "C23"
$ws.Range("A16").Value = "C24"
$ws.Range("A17").Value = "C25"
$ws.Range("A18").Value = "C26"
$ws.Range("A19").Value = "C27"
$ws.Range("A20").Value = "C28"
$ws.Range("A21").Value = "C29"
$ws.Range("A22").Value = "C30"
$ws.Range("A23").Value = "C31-32"
$ws.Range("A24").Value = "C33"
$ws.Range("A25").Value = "D35"
$ws.Range("A26").Value = "E36"
$ws.Range("A27").Value = "E37-39"
$ws.Range("A28").Value = "F41-43"
$ws.Range("A29").Value = "G45"
$ws.Range("A30").Value = "G46"
$ws.Range("A31").Value = "G47"
$ws.Range("A32").Value = "H49"
$ws.Range("A33").Value = "H50"
$ws.Range("A34").Value = "H51"
$ws.Range("A35").Value = "H52"
$ws.Range("A36").Value = "H53"
$ws.Range("A37").Value = "I55-56"
$ws.Range("A38").Value = "J58"
$ws.Range("A39").Value = "J59-60"
$ws.Range("A40").Value = "J61"
$ws.Range("A41").Value = "J62-63"
$ws.Range("A42").Value = "K64"
$ws.Range("A43").Value = "K65"
$ws.Range("A44").Value = "K66"
$ws.Range("A45").Value = "L68"
$ws.Range("A46").Value = "M69-70"
$ws.Range("A47").Value = "M71"
$ws.Range("A48").Value = "M72"
$ws.Range("A49").Value = "M73"
$ws.Range("A50").Value = "M74-75"
$ws.Range("A51").Value = "N77"
$ws.Range("A52").Value = "N78"
$ws.Range("A53").Value = "N79"
$ws.Range("A54").Value = "N80-82"
$ws.Range("A55").Value = "O84"
$ws.Range("A56").Value = "P85"
$ws.Range("A57").Value = "Q86"
$ws.Range("A58").Value = "Q87-88"
$ws.Range("A59").Value = "R90-92"
$ws.Range("A60").Value = "R93"
$ws.Range("A61").Value = "S94"
$ws.Range("A62").Value = "S95"
$ws.Range("A63").Value = "S96"
$ws.Range("A64").Value = "T97-98"

# Add 16 blank formatted rows below the data (rows 66-81)
for ($r = 66; $r -le 81; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Move the selection to A66, matching the post-edit cursor position
[void]$ws.Range("A66").Select()
